$wb = $excel.ActiveWorkbook

# --- Sheet 1 ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(159, 1).Value = 45945.49385416666
$ws.Cells.Item(159, 1).NumberFormat = $ws.Cells.Item(158, 1).NumberFormat
$ws.Cells.Item(159, 2).Value = "0x01,0x7c"
$ws.Cells.Item(159, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(159, 4).Value = "0x00,0xB8"
$ws.Cells.Item(159, 5).Value = "0xf"
$ws.Cells.Item(159, 6).Value = 380
$ws.Cells.Item(159, 7).Value = 759863127514711000000000.0
$ws.Cells.Item(159, 8).Value = 196
$ws.Cells.Item(159, 9).Value = 15

$ws.Cells.Item(160, 1).Value = 45946.49204861111
$ws.Cells.Item(160, 1).NumberFormat = $ws.Cells.Item(158, 1).NumberFormat
$ws.Cells.Item(160, 2).Value = "0x01,0x7c"
$ws.Cells.Item(160, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item(160, 4).Value = "0x00,0xB8"
$ws.Cells.Item(160, 5).Value = "0xf"
$ws.Cells.Item(160, 6).Value = 380
$ws.Cells.Item(160, 7).Value = 759863127514711000000000.0
$ws.Cells.Item(160, 8).Value = 192
$ws.Cells.Item(160, 9).Value = 15

# --- Sheet 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(159, 1).Value = 45945.49385416666
$ws.Cells.Item(159, 1).NumberFormat = $ws.Cells.Item(158, 1).NumberFormat
$ws.Cells.Item(159, 2).Value = "0x01,0x90"
$ws.Cells.Item(159, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(159, 4).Value = "0x00,0xBC"
$ws.Cells.Item(159, 5).Value = "0xe"
$ws.Cells.Item(159, 6).Value = 400
$ws.Cells.Item(159, 7).Value = 568432987514711000000000.0
$ws.Cells.Item(159, 8).Value = 196
$ws.Cells.Item(159, 9).Value = 14

$ws.Cells.Item(160, 1).Value = 45946.49204861111
$ws.Cells.Item(160, 1).NumberFormat = $ws.Cells.Item(158, 1).NumberFormat
$ws.Cells.Item(160, 2).Value = "0x01,0x90"
$ws.Cells.Item(160, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item(160, 4).Value = "0x00,0xB8"
$ws.Cells.Item(160, 5).Value = "0xe"
$ws.Cells.Item(160, 6).Value = 400
$ws.Cells.Item(160, 7).Value = 568432987514711000000000.0
$ws.Cells.Item(160, 8).Value = 196
$ws.Cells.Item(160, 9).Value = 14

# --- Sheet 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(159, 1).Value = 45945.49385416666
$ws.Cells.Item(159, 1).NumberFormat = $ws.Cells.Item(158, 1).NumberFormat
$ws.Cells.Item(159, 2).Value = "0x00,0x6e"
$ws.Cells.Item(159, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(159, 4).Value = "0x00,0x50"
$ws.Cells.Item(159, 5).Value = "0x3"
$ws.Cells.Item(159, 6).Value = 110
$ws.Cells.Item(159, 7).Value = 568631262647114000000000.0
$ws.Cells.Item(159, 8).Value = 80
$ws.Cells.Item(159, 9).Value = 3

$ws.Cells.Item(160, 1).Value = 45946.49204861111
$ws.Cells.Item(160, 1).NumberFormat = $ws.Cells.Item(158, 1).NumberFormat
$ws.Cells.Item(160, 2).Value = "0x00,0x6e"
$ws.Cells.Item(160, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item(160, 4).Value = "0x00,0x4F"
$ws.Cells.Item(160, 5).Value = "0x3"
$ws.Cells.Item(160, 6).Value = 110
$ws.Cells.Item(160, 7).Value = 568631262647114000000000.0
$ws.Cells.Item(160, 8).Value = 79
$ws.Cells.Item(160, 9).Value = 3

# --- Sheet 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(159, 1).Value = 45945.49385416666
$ws.Cells.Item(159, 1).NumberFormat = $ws.Cells.Item(158, 1).NumberFormat
$ws.Cells.Item(159, 2).Value = "0x00,0x6e"
$ws.Cells.Item(159, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(159, 4).Value = "0x00,0x4C"
$ws.Cells.Item(159, 5).Value = "0x3"
$ws.Cells.Item(159, 6).Value = 110
$ws.Cells.Item(159, 7).Value = 985046333984776000000000.0
$ws.Cells.Item(159, 8).Value = 76
$ws.Cells.Item(159, 9).Value = 3

$ws.Cells.Item(160, 1).Value = 45946.49204861111
$ws.Cells.Item(160, 1).NumberFormat = $ws.Cells.Item(158, 1).NumberFormat
$ws.Cells.Item(160, 2).Value = "0x00,0x6e"
$ws.Cells.Item(160, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item(160, 4).Value = "0x00,0x4C"
$ws.Cells.Item(160, 5).Value = "0x3"
$ws.Cells.Item(160, 6).Value = 110
$ws.Cells.Item(160, 7).Value = 985046333984776000000000.0
$ws.Cells.Item(160, 8).Value = 76
$ws.Cells.Item(160, 9).Value = 3
